$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: 'Morning Glass of Ether' | 'Ether'
$ws.Range("H15").Value = 5175.6025
$ws.Range("I15").Value = 5175.6025
$ws.Range("K15").Value = 15526.8075
$ws.Range("M15").Value = -15357.8075

# Row 125: 'Body over Mind' | 'Grade 5 Dexterity Alkahest'
$ws.Range("H125").Value = 483.89474
$ws.Range("J125").Value = 301.7647
$ws.Range("L125").Value = 2715.8823
$ws.Range("N125").Value = -7635.8823

# Row 126: 'Rebuilding to Code' | 'Saigaskin Codex'
$ws.Range("H126").Value = 46769.332
$ws.Range("J126").Value = 46769.332
$ws.Range("L126").Value = 46769.332
$ws.Range("N126").Value = -56649.332

# Row 128: 'Nearly There' | 'Kumbhiraskin Grimoire'
$ws.Range("H128").Value = 42685.832
$ws.Range("J128").Value = 42685.832
$ws.Range("L128").Value = 42685.832
$ws.Range("N128").Value = -52645.832

# Row 130: 'Technically Still Magic' | 'Ophiotauroskin Magitek Codex'
$ws.Range("H130").Value = 44003.2
$ws.Range("J130").Value = 44003.2
$ws.Range("L130").Value = 44003.2
$ws.Range("N130").Value = -54043.2

$ws = $wb.Worksheets.Item("ARM")
# Row 123: 'The Armoire Is Open' | 'High Durium Armguards of Maiming'
$ws.Range("H123").Value = 35614.5
$ws.Range("J123").Value = 35614.5
$ws.Range("L123").Value = 35614.5
$ws.Range("N123").Value = -45414.5

# Row 128: 'Heading toward Bankruptcy' | 'Manganese Helm of the Falling Dragon'
$ws.Range("H128").Value = 47090.332
$ws.Range("J128").Value = 47090.332
$ws.Range("L128").Value = 47090.332
$ws.Range("N128").Value = -57050.332

# Row 130: 'A Gift of Gloves' | 'Chondrite Gloves of Casting'
$ws.Range("H130").Value = 44996
$ws.Range("J130").Value = 44996
$ws.Range("L130").Value = 44996
$ws.Range("N130").Value = -55036

# Row 131: 'Additions to the Armoire' | 'Chondrite Top of Maiming'
$ws.Range("H131").Value = 49607
$ws.Range("J131").Value = 49607
$ws.Range("L131").Value = 49607
$ws.Range("N131").Value = -59687

$ws = $wb.Worksheets.Item("BSM")
# Row 20: 'Smelt and Dealt' | 'Iron Ingot'
$ws.Range("H20").Value = 1685.7368
$ws.Range("I20").Value = 1221.8182
$ws.Range("J20").Value = 2323.625
$ws.Range("K20").Value = 1221.8182
$ws.Range("L20").Value = 2323.625
$ws.Range("M20").Value = -974.8181999999999
$ws.Range("N20").Value = -2817.625

# Row 86: 'Through Thick and Thin' | 'Adamantite Nugget'
$ws.Range("H86").Value = 2493.2307
$ws.Range("I86").Value = 2572.1428
$ws.Range("J86").Value = 2401.1667
$ws.Range("K86").Value = 2572.1428
$ws.Range("L86").Value = 2401.1667
$ws.Range("M86").Value = -1449.1428
$ws.Range("N86").Value = -4647.1667

# Row 89: 'Piercing Eyes Deserve Piercing Shafts (L)' | 'Adamantite Nugget'
$ws.Range("H89").Value = 2493.2307
$ws.Range("I89").Value = 2572.1428
$ws.Range("J89").Value = 2401.1667
$ws.Range("K89").Value = 12860.714
$ws.Range("L89").Value = 12005.8335
$ws.Range("M89").Value = -7244.714
$ws.Range("N89").Value = -23237.8335

# Row 105: 'Ingot to Wing It' | 'Molybdenum Ingot'
$ws.Range("H105").Value = 2204.8667
$ws.Range("I105").Value = 1972.2222
$ws.Range("K105").Value = 1972.2222
$ws.Range("M105").Value = -225.2221999999999

# Row 124: 'History of the Hrothgar' | 'High Durium Bayonet'
$ws.Range("H124").Value = 49996
$ws.Range("J124").Value = 49996
$ws.Range("L124").Value = 49996
$ws.Range("N124").Value = -59816

# Row 125: 'Archon of His Eye' | 'High Durium Knives'
$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620

# Row 126: 'Records of the Republic' | 'Bismuth War Scythe'
$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656

# Row 130: 'Annals of the Empire I' | 'Chondrite Magitek Axe'
$ws.Range("H130").Value = 45077
$ws.Range("J130").Value = 45077
$ws.Range("L130").Value = 45077
$ws.Range("N130").Value = -55117

# Row 134: 'Ruthenium Supremium' | 'Ruthenium Ingot'
$ws.Range("H134").Value = 4692.4253
$ws.Range("I134").Value = 4749.6665
$ws.Range("J134").Value = 4678.8687
$ws.Range("K134").Value = 14248.9995
$ws.Range("L134").Value = 14036.6061
$ws.Range("M134").Value = -11713.9995
$ws.Range("N134").Value = -19106.6061

$ws = $wb.Worksheets.Item("CRP")
# Row 20: 'Re-crating the Scene' | 'Iron Spear'
$ws.Range("H20").Value = 40794.668
$ws.Range("J20").Value = 40794.668
$ws.Range("L20").Value = 40794.668
$ws.Range("N20").Value = -41266.668

# Row 30: 'Polearms Aplenty' | 'Iron Spear'
$ws.Range("H30").Value = 40794.668
$ws.Range("J30").Value = 40794.668
$ws.Range("L30").Value = 40794.668
$ws.Range("N30").Value = -40976.668

# Row 31: 'Wall Not Found' | 'Walnut Lumber'
$ws.Range("H31").Value = 146084.36
$ws.Range("I31").Value = 1952.2941
$ws.Range("J31").Value = 175965.4
$ws.Range("K31").Value = 1952.2941
$ws.Range("L31").Value = 175965.4
$ws.Range("M31").Value = -1657.2941
$ws.Range("N31").Value = -176555.4

# Row 34: 'Armoires of the Rich and Famous' | 'Walnut Lumber'
$ws.Range("H34").Value = 146084.36
$ws.Range("I34").Value = 1952.2941
$ws.Range("J34").Value = 175965.4
$ws.Range("K34").Value = 1952.2941
$ws.Range("L34").Value = 175965.4
$ws.Range("M34").Value = -1750.2941
$ws.Range("N34").Value = -176369.4

# Row 128: 'An A-prop-riate Request' | 'Ironwood Spear'
$ws.Range("H128").Value = 40794.668
$ws.Range("J128").Value = 40794.668
$ws.Range("L128").Value = 40794.668
$ws.Range("N128").Value = -50754.668

$ws = $wb.Worksheets.Item("CUL")
# Row 97: 'The Frier Never Lies' | 'Cottonseed Oil'
$ws.Range("H97").Value = 1499
$ws.Range("J97").Value = 1499
$ws.Range("L97").Value = 4497
$ws.Range("N97").Value = -5489

# Row 117: 'A Good Omen' | 'Peppered Popotoes'
$ws.Range("H117").Value = 1000
$ws.Range("I117").Value = 877.7778
$ws.Range("J117").Value = 1366.6666
$ws.Range("K117").Value = 2633.3334
$ws.Range("L117").Value = 4099.9998
$ws.Range("M117").Value = 808.6666
$ws.Range("N117").Value = -10983.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 58: 'The Big Red' | 'Red Coral Necklace'
$ws.Range("H58").Value = 24975
$ws.Range("J58").Value = 24975
$ws.Range("L58").Value = 24975
$ws.Range("N58").Value = -25529

# Row 110: 'Slimming Down' | 'Stonegold Rapier'
$ws.Range("H110").Value = 43465.332
$ws.Range("J110").Value = 43465.332
$ws.Range("L110").Value = 43465.332
$ws.Range("N110").Value = -51645.332

# Row 130: 'Planisphere to Paper' | 'Chondrite Magitek Planisphere'
$ws.Range("H130").Value = 45784
$ws.Range("J130").Value = 45784
$ws.Range("L130").Value = 45784
$ws.Range("N130").Value = -55824

$ws = $wb.Worksheets.Item("LTW")
# Row 36: 'Campaign in the Membrane' | 'Toadskin Jacket'
$ws.Range("H36").Value = 45807.332
$ws.Range("J36").Value = 45807.332
$ws.Range("L36").Value = 45807.332
$ws.Range("N36").Value = -46931.332

# Row 57: 'Too Hot to Handle' | 'Raptorskin Wristbands'
$ws.Range("H57").Value = 37925.668
$ws.Range("J57").Value = 37925.668
$ws.Range("L57").Value = 37925.668
$ws.Range("N57").Value = -39057.668

# Row 122: 'Hell on Leather' | 'Gaja Leather'
$ws.Range("H122").Value = 2076.4
$ws.Range("I122").Value = 1928.5333
$ws.Range("K122").Value = 5785.5999
$ws.Range("M122").Value = -3335.5999

# Row 127: 'Loyal Turncoat' | 'Saigaskin Coat of Fending'
$ws.Range("H127").Value = 50608.832
$ws.Range("J127").Value = 50608.832
$ws.Range("L127").Value = 50608.832
$ws.Range("N127").Value = -60528.832

# Row 130: 'Generous Soles' | 'Ophiotauroskin Boots of Healing'
$ws.Range("H130").Value = 37996
$ws.Range("J130").Value = 37996
$ws.Range("L130").Value = 37996
$ws.Range("N130").Value = -48036

# Row 136: "Respect for Br'aax" | "Br'aax Leather"
$ws.Range("H136").Value = 2109.4
$ws.Range("I136").Value = 1446
$ws.Range("J136").Value = 4099.6
$ws.Range("K136").Value = 4338
$ws.Range("L136").Value = 12298.8
$ws.Range("M136").Value = -1788
$ws.Range("N136").Value = -17398.8

$ws = $wb.Worksheets.Item("WVR")
# Row 120: 'A Turban for the Ages' | 'Dwarven Cotton Turban of Scouting'
$ws.Range("H120").Value = 40206
$ws.Range("J120").Value = 40206
$ws.Range("L120").Value = 40206
$ws.Range("N120").Value = -49882

# Row 122: 'Heavy Armoire' | 'Dark Hempen Cloth'
$ws.Range("H122").Value = 2599076
$ws.Range("I122").Value = 4083048
$ws.Range("J122").Value = 2125
$ws.Range("K122").Value = 12249144
$ws.Range("L122").Value = 6375
$ws.Range("M122").Value = -12246694
$ws.Range("N122").Value = -11275

# Row 124: 'Hot Heads' | 'Almasty Serge Hat of Casting'
$ws.Range("H124").Value = 24357.25
$ws.Range("J124").Value = 24357.25
$ws.Range("L124").Value = 24357.25
$ws.Range("N124").Value = -34177.25

# Row 128: 'Lightening Up' | 'Scarlet Moko Gaskins of the Rising Dragon'
$ws.Range("H128").Value = 49715
$ws.Range("J128").Value = 49715
$ws.Range("L128").Value = 49715
$ws.Range("N128").Value = -59675

# Row 136: 'Weaving the Envelope' | 'Sarcenet Cloth'
$ws.Range("H136").Value = 18563.492
$ws.Range("I136").Value = 30981.94
$ws.Range("J136").Value = 1488.125
$ws.Range("K136").Value = 92945.81999999999
$ws.Range("L136").Value = 4464.375
$ws.Range("M136").Value = -90395.81999999999
$ws.Range("N136").Value = -9564.375
